$wb = $excel.ActiveWorkbook

# Sheet ALC row 15
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H15").Value = 632.88
$ws.Range("I15").Value = 632.88
$ws.Range("J15").Value = 0
$ws.Range("K15").Value = 1898.64
$ws.Range("L15").Value = 0
$ws.Range("M15").Value = -1729.64

# Sheet ALC row 40
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1339.9878
$ws.Range("I40").Value = 1247.28
$ws.Range("J40").Value = 2333.2856
$ws.Range("K40").Value = 1247.28
$ws.Range("L40").Value = 2333.2856
$ws.Range("M40").Value = -1072.28
$ws.Range("N40").Value = -2683.2856

# Sheet ALC row 47
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H47").Value = 3000
$ws.Range("I47").Value = 3000
$ws.Range("J47").Value = 0
$ws.Range("K47").Value = 3000
$ws.Range("L47").Value = 0
$ws.Range("M47").Value = -2028

# Sheet ALC row 129
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H129").Value = 982.4375
$ws.Range("I129").Value = 385.7143
$ws.Range("J129").Value = 1446.5555
$ws.Range("K129").Value = 1157.1429
$ws.Range("L129").Value = 4339.666499999999
$ws.Range("M129").Value = 3842.8571
$ws.Range("N129").Value = -14339.6665

# Sheet ALC row 130
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H130").Value = 35800
$ws.Range("I130").Value = 0
$ws.Range("J130").Value = 35800
$ws.Range("K130").Value = 0
$ws.Range("L130").Value = 35800
$ws.Range("N130").Value = -45840

# Sheet ALC row 132
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 2496.3247
$ws.Range("I132").Value = 1171.8572
$ws.Range("J132").Value = 8456.429
$ws.Range("K132").Value = 3515.5716
$ws.Range("L132").Value = 25369.287
$ws.Range("M132").Value = -985.5715999999998
$ws.Range("N132").Value = -30429.287

# Sheet ALC row 133
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H133").Value = 28010.889
$ws.Range("I133").Value = 0
$ws.Range("J133").Value = 28010.889
$ws.Range("K133").Value = 0
$ws.Range("L133").Value = 28010.889
$ws.Range("N133").Value = -38130.889

# Sheet ALC row 138
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2289.3157
$ws.Range("I138").Value = 1319.4
$ws.Range("J138").Value = 2635.7144
$ws.Range("K138").Value = 3958.2
$ws.Range("L138").Value = 7907.1432
$ws.Range("M138").Value = 1181.8
$ws.Range("N138").Value = -18187.1432

# Sheet ARM row 32
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5991.282
$ws.Range("I32").Value = 4466.5415
$ws.Range("J32").Value = 14436
$ws.Range("K32").Value = 4466.5415
$ws.Range("L32").Value = 14436
$ws.Range("M32").Value = -4179.5415
$ws.Range("N32").Value = -15010

# Sheet ARM row 123
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H123").Value = 24999
$ws.Range("I123").Value = 0
$ws.Range("J123").Value = 24999
$ws.Range("K123").Value = 0
$ws.Range("L123").Value = 24999
$ws.Range("N123").Value = -34799

# Sheet ARM row 124
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 19445.8
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 19445.8
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 19445.8
$ws.Range("N124").Value = -29265.8

# Sheet ARM row 125
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H125").Value = 0
$ws.Range("I125").Value = 0
$ws.Range("J125").Value = 0
$ws.Range("K125").Value = 0
$ws.Range("L125").Value = 0
$ws.Range("N125").ClearContents()

# Sheet BSM row 94
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H94").Value = 1578.9
$ws.Range("I94").Value = 1071.4286
$ws.Range("J94").Value = 1852.1538
$ws.Range("K94").Value = 1071.4286
$ws.Range("L94").Value = 1852.1538
$ws.Range("M94").Value = -620.4286
$ws.Range("N94").Value = -2754.1538

# Sheet BSM row 124
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H124").Value = 0
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 0
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 0
$ws.Range("N124").ClearContents()

# Sheet CRP row 42
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H42").Value = 10000
$ws.Range("I42").Value = 10000
$ws.Range("J42").Value = 0
$ws.Range("K42").Value = 10000
$ws.Range("L42").Value = 0
$ws.Range("M42").Value = -9407

# Sheet CRP row 110
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H110").Value = 26000
$ws.Range("I110").Value = 0
$ws.Range("J110").Value = 26000
$ws.Range("K110").Value = 0
$ws.Range("L110").Value = 26000
$ws.Range("N110").Value = -34180

# Sheet CRP row 111
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H111").Value = 29435
$ws.Range("I111").Value = 0
$ws.Range("J111").Value = 29435
$ws.Range("K111").Value = 0
$ws.Range("L111").Value = 29435
$ws.Range("N111").Value = -37615

# Sheet CRP row 114
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H114").Value = 29666.666
$ws.Range("I114").Value = 0
$ws.Range("J114").Value = 29666.666
$ws.Range("K114").Value = 0
$ws.Range("L114").Value = 29666.666
$ws.Range("N114").Value = -38344.666

# Sheet CRP row 116
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H116").Value = 39000
$ws.Range("I116").Value = 0
$ws.Range("J116").Value = 39000
$ws.Range("K116").Value = 0
$ws.Range("L116").Value = 39000
$ws.Range("N116").Value = -48178

# Sheet CRP row 122
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H122").Value = 1226.4286
$ws.Range("I122").Value = 759.5454999999999
$ws.Range("J122").Value = 1740
$ws.Range("K122").Value = 2278.6365
$ws.Range("L122").Value = 5220
$ws.Range("M122").Value = 171.3635000000004
$ws.Range("N122").Value = -10120

# Sheet CRP row 124
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H124").Value = 17927.143
$ws.Range("I124").Value = 0
$ws.Range("J124").Value = 17927.143
$ws.Range("K124").Value = 0
$ws.Range("L124").Value = 17927.143
$ws.Range("N124").Value = -22837.143

# Sheet CUL row 39
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H39").Value = 2541.1765
$ws.Range("I39").Value = 600
$ws.Range("J39").Value = 2800
$ws.Range("K39").Value = 1800
$ws.Range("L39").Value = 8400
$ws.Range("M39").Value = -1506
$ws.Range("N39").Value = -8988

# Sheet CUL row 42
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H42").Value = 3033.3333
$ws.Range("I42").Value = 0
$ws.Range("J42").Value = 3033.3333
$ws.Range("K42").Value = 0
$ws.Range("L42").Value = 9099.999899999999
$ws.Range("N42").Value = -10167.9999

# Sheet GSM row 97
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1216.3636
$ws.Range("I97").Value = 935.3077
$ws.Range("J97").Value = 1622.3334
$ws.Range("K97").Value = 935.3077
$ws.Range("L97").Value = 1622.3334
$ws.Range("M97").Value = -439.3077
$ws.Range("N97").Value = -2614.3334

# Sheet GSM row 102
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 2949.16
$ws.Range("I102").Value = 3041.261
$ws.Range("J102").Value = 1890
$ws.Range("K102").Value = 3041.261
$ws.Range("L102").Value = 1890
$ws.Range("M102").Value = -1419.261
$ws.Range("N102").Value = -5134

# Sheet GSM row 122
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1062.9286
$ws.Range("I122").Value = 940.0833
$ws.Range("J122").Value = 1800
$ws.Range("K122").Value = 2820.2499
$ws.Range("L122").Value = 5400
$ws.Range("M122").Value = -370.2498999999998
$ws.Range("N122").Value = -10300

# Sheet GSM row 132
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 5680.5835
$ws.Range("I132").Value = 7117.6665
$ws.Range("J132").Value = 4243.5
$ws.Range("K132").Value = 21352.9995
$ws.Range("L132").Value = 12730.5
$ws.Range("M132").Value = -18822.9995
$ws.Range("N132").Value = -17790.5

# Sheet LTW row 22
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 167.04
$ws.Range("I22").Value = 165.3
$ws.Range("J22").Value = 174
$ws.Range("K22").Value = 165.3
$ws.Range("L22").Value = 174
$ws.Range("M22").Value = 129.7
$ws.Range("N22").Value = -764

# Sheet LTW row 27
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 167.04
$ws.Range("I27").Value = 165.3
$ws.Range("J27").Value = 174
$ws.Range("K27").Value = 165.3
$ws.Range("L27").Value = 174
$ws.Range("M27").Value = -58.30000000000001
$ws.Range("N27").Value = -388

# Sheet LTW row 40
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 2813.0667
$ws.Range("I40").Value = 2717.818
$ws.Range("J40").Value = 3075
$ws.Range("K40").Value = 2717.818
$ws.Range("L40").Value = 3075
$ws.Range("M40").Value = -2581.818
$ws.Range("N40").Value = -3347

# Sheet LTW row 41
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H41").Value = 12900
$ws.Range("I41").Value = 0
$ws.Range("J41").Value = 12900
$ws.Range("K41").Value = 0
$ws.Range("L41").Value = 12900
$ws.Range("N41").Value = -13776

# Sheet LTW row 46
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1125
$ws.Range("I46").Value = 716.6667
$ws.Range("J46").Value = 1533.3334
$ws.Range("K46").Value = 716.6667
$ws.Range("L46").Value = 1533.3334
$ws.Range("M46").Value = -528.6667
$ws.Range("N46").Value = -1909.3334

# Sheet LTW row 48
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H48").Value = 5805.5
$ws.Range("I48").Value = 920.5
$ws.Range("J48").Value = 8248
$ws.Range("K48").Value = 920.5
$ws.Range("L48").Value = 8248
$ws.Range("M48").Value = -259.5
$ws.Range("N48").Value = -9570

# Sheet LTW row 100
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 40004344
$ws.Range("I100").Value = 7004.3
$ws.Range("J100").Value = 66669236
$ws.Range("K100").Value = 7004.3
$ws.Range("L100").Value = 66669236
$ws.Range("M100").Value = -6463.3
$ws.Range("N100").Value = -66670318

# Sheet LTW row 122
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 2433.6072
$ws.Range("I122").Value = 2012.6
$ws.Range("J122").Value = 2919.3845
$ws.Range("K122").Value = 6037.799999999999
$ws.Range("L122").Value = 8758.1535
$ws.Range("M122").Value = -3587.799999999999
$ws.Range("N122").Value = -13658.1535

# Sheet WVR row 96
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 7485.75
$ws.Range("I96").Value = 0
$ws.Range("J96").Value = 7485.75
$ws.Range("K96").Value = 0
$ws.Range("L96").Value = 7485.75
$ws.Range("M96").ClearContents()
$ws.Range("N96").Value = -10231.75

# Sheet WVR row 122
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 49033.332
$ws.Range("I122").Value = 59964.707
$ws.Range("J122").Value = 2575
$ws.Range("K122").Value = 179894.121
$ws.Range("L122").Value = 7725
$ws.Range("M122").Value = -177444.121
$ws.Range("N122").Value = -12625

# Sheet WVR row 132
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 26787.463
$ws.Range("I132").Value = 60144.59
$ws.Range("J132").Value = 3159.5
$ws.Range("K132").Value = 180433.77
$ws.Range("L132").Value = 9478.5
$ws.Range("M132").Value = -177903.77
$ws.Range("N132").Value = -14538.5

# Sheet WVR row 136
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 29631940
$ws.Range("I136").Value = 34484590
$ws.Range("J136").Value = 20836508
$ws.Range("K136").Value = 103453770
$ws.Range("L136").Value = 62509524
$ws.Range("M136").Value = -103451220
$ws.Range("N136").Value = -62514624
